$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2501163
$ws.Range("J17").Value = 2501163
$ws.Range("L17").Value = 7503489
$ws.Range("N17").Value = -7503825
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H33").Value = 1501.5
$ws.Range("I33").Value = 1001.6667
$ws.Range("K33").Value = 1001.6667
$ws.Range("M33").Value = -772.6667
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("H137").Value = 1424.3636
$ws.Range("J137").Value = 1672.5
$ws.Range("L137").Value = 5017.5
$ws.Range("N137").Value = -10117.5
$ws.Range("H139").Value = 55000
$ws.Range("J139").Value = 55000
$ws.Range("L139").Value = 55000
$ws.Range("N139").Value = -65280
$ws.Range("H141").Value = 3532.6365
$ws.Range("I141").Value = 3490
$ws.Range("J141").Value = 3646.3333
$ws.Range("K141").Value = 10470
$ws.Range("L141").Value = 10938.9999
$ws.Range("M141").Value = -5290
$ws.Range("N141").Value = -21298.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 3733.3333
$ws.Range("I26").Value = 3733.3333
$ws.Range("K26").Value = 3733.3333
$ws.Range("M26").Value = -3403.3333
$ws.Range("H61").Value = 9691.637000000001
$ws.Range("I61").Value = 9419.75
$ws.Range("K61").Value = 9419.75
$ws.Range("M61").Value = -9207.75
$ws.Range("H97").Value = 903.94446
$ws.Range("J97").Value = 2974.5
$ws.Range("L97").Value = 2974.5
$ws.Range("N97").Value = -3966.5
$ws.Range("H106").Value = 40365
$ws.Range("J106").Value = 40365
$ws.Range("L106").Value = 40365
$ws.Range("N106").Value = -42889
$ws.Range("H132").Value = 3941.5881
$ws.Range("I132").Value = 3323.093
$ws.Range("J132").Value = 7266
$ws.Range("K132").Value = 9969.278999999999
$ws.Range("L132").Value = 21798
$ws.Range("M132").Value = -7439.278999999999
$ws.Range("N132").Value = -26858
$ws.Range("H136").Value = 9691.637000000001
$ws.Range("I136").Value = 9419.75
$ws.Range("K136").Value = 28259.25
$ws.Range("M136").Value = -25709.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 31028.334
$ws.Range("J35").Value = 32074
$ws.Range("L35").Value = 32074
$ws.Range("N35").Value = -32694
$ws.Range("H107").Value = 1608.5555
$ws.Range("I107").Value = 1190.4
$ws.Range("K107").Value = 1190.4
$ws.Range("M107").Value = 729.5999999999999
$ws.Range("H135").Value = 119999.5
$ws.Range("J135").Value = 119999.5
$ws.Range("L135").Value = 119999.5
$ws.Range("N135").Value = -130139.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 306
$ws.Range("I19").Value = 370.625
$ws.Range("J19").Value = 47.5
$ws.Range("K19").Value = 370.625
$ws.Range("L19").Value = 47.5
$ws.Range("M19").Value = -200.625
$ws.Range("N19").Value = -387.5
$ws.Range("H24").Value = 306
$ws.Range("I24").Value = 370.625
$ws.Range("J24").Value = 47.5
$ws.Range("K24").Value = 370.625
$ws.Range("L24").Value = 47.5
$ws.Range("M24").Value = -200.625
$ws.Range("N24").Value = -387.5
$ws.Range("H58").Value = 4662.9688
$ws.Range("I58").Value = 2959.762
$ws.Range("J58").Value = 7914.5454
$ws.Range("K58").Value = 2959.762
$ws.Range("L58").Value = 7914.5454
$ws.Range("M58").Value = -2756.762
$ws.Range("N58").Value = -8320.545399999999
$ws.Range("I59").Value = 40000
$ws.Range("K59").Value = 40000
$ws.Range("M59").Value = -38855
$ws.Range("H68").Value = 41641.332
$ws.Range("J68").Value = 41641.332
$ws.Range("L68").Value = 41641.332
$ws.Range("N68").Value = -43139.332
$ws.Range("H71").Value = 41641.332
$ws.Range("J71").Value = 41641.332
$ws.Range("L71").Value = 124923.996
$ws.Range("N71").Value = -132411.996
$ws.Range("H74").Value = 40995.6
$ws.Range("J74").Value = 40995.6
$ws.Range("L74").Value = 40995.6
$ws.Range("N74").Value = -42743.6
$ws.Range("H77").Value = 40995.6
$ws.Range("J77").Value = 40995.6
$ws.Range("L77").Value = 122986.8
$ws.Range("N77").Value = -131722.8
$ws.Range("H133").Value = 51343.285
$ws.Range("J133").Value = 54080.6
$ws.Range("L133").Value = 54080.6
$ws.Range("N133").Value = -59140.6
$ws.Range("H134").Value = 5959.905
$ws.Range("I134").Value = 5008.316
$ws.Range("J134").Value = 15000
$ws.Range("K134").Value = 15024.948
$ws.Range("L134").Value = 45000
$ws.Range("M134").Value = -12489.948
$ws.Range("N134").Value = -50070
$ws.Range("H136").Value = 4662.9688
$ws.Range("I136").Value = 2959.762
$ws.Range("J136").Value = 7914.5454
$ws.Range("K136").Value = 8879.286
$ws.Range("L136").Value = 23743.6362
$ws.Range("M136").Value = -6329.286
$ws.Range("N136").Value = -28843.6362

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 40257024
$ws.Range("I4").Value = 36523090
$ws.Range("K4").Value = 109569270
$ws.Range("M4").Value = -109569158
$ws.Range("H5").Value = 817.6842
$ws.Range("I5").Value = 431.35715
$ws.Range("K5").Value = 1294.07145
$ws.Range("M5").Value = -1182.07145
$ws.Range("H133").Value = 2750
$ws.Range("I133").Value = 2750
$ws.Range("K133").Value = 8250
$ws.Range("M133").Value = -3190
$ws.Range("H135").Value = 817.6842
$ws.Range("I135").Value = 431.35715
$ws.Range("K135").Value = 3882.21435
$ws.Range("M135").Value = -1347.21435
$ws.Range("H137").Value = 1314.25
$ws.Range("I137").Value = 1216.3572
$ws.Range("J137").Value = 1999.5
$ws.Range("K137").Value = 3649.0716
$ws.Range("L137").Value = 5998.5
$ws.Range("M137").Value = 1450.9284
$ws.Range("N137").Value = -16198.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 38806.855
$ws.Range("I46").Value = 32333.334
$ws.Range("J46").Value = 43662
$ws.Range("K46").Value = 32333.334
$ws.Range("L46").Value = 43662
$ws.Range("M46").Value = -32177.334
$ws.Range("N46").Value = -43974
$ws.Range("H132").Value = 9833
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 9833
$ws.Range("K132").Value = 0
$ws.Range("L132").ClearContents()
$ws.Range("M132").Value = 29499
$ws.Range("N132").Value = -34559

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1910
$ws.Range("I16").Value = 1733.3334
$ws.Range("K16").Value = 1733.3334
$ws.Range("M16").Value = -1563.3334
$ws.Range("H32").Value = 12500
$ws.Range("I32").Value = 12500
$ws.Range("K32").Value = 12500
$ws.Range("M32").Value = -12183
$ws.Range("H136").Value = 7282.1177
$ws.Range("I136").Value = 7622.1113
$ws.Range("J136").Value = 6899.625
$ws.Range("K136").Value = 22866.3339
$ws.Range("L136").Value = 20698.875
$ws.Range("M136").Value = -20316.3339
$ws.Range("N136").Value = -25798.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2429.0476
$ws.Range("I122").Value = 1494
$ws.Range("K122").Value = 4482
$ws.Range("M122").Value = -2032
$ws.Range("H126").Value = 3993.7727
$ws.Range("I126").Value = 3727
$ws.Range("K126").Value = 11181
$ws.Range("M126").Value = -8711
$ws.Range("H132").Value = 2901.9473
$ws.Range("J132").Value = 6135.8
$ws.Range("L132").Value = 18407.4
$ws.Range("N132").Value = -23467.4
